# Apply the updated DSD metadata values to the sheet.
# The whole 4x8 block (A1:H4) is rewritten with corrected/reordered
# header labels, measure/dimension identifiers, datatypes and URI columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @("Tipo hogar 1", "Tipo de hogar 2", "Número hogares", "Provincia código", "Aragón", "Municipio código", "Provincia nombre", "Municipio nombre"),
    @("iaest-measure:tipo-hogar-1", "iaest-measure:tipo-de-hogar-2", "iaest-measure:numero-hogares", "null", "sdmx-dimension:refArea", "null", "sdmx-dimension:refArea", "sdmx-dimension:refArea"),
    @("medida", "medida", "medida", "null", "dim", "null", "dim", "dim"),
    @("xsd:string", "xsd:string", "xsd:int", "null", "URI-Comunidad", "null", "URI-Provincia", "URI-Municipio")
)

for ($r = 0; $r -lt 4; $r++) {
    for ($c = 0; $c -lt 8; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $values[$r][$c]
    }
}

$wb.Save()
